$p = $ppt.ActivePresentation

# --- Slide 1: remove the empty "Subtitle 2" placeholder shape ---
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "Subtitle 2") {
        $shp.Delete()
    }
}

# --- Slide 13: merge the two runs "Luiz Antonio Pereira Silva, " and
#     "Joaquim José Cintra Maia Honório" into a single run/text ---
$s13 = $p.Slides.Item(13)
$tbl = $s13.Shapes.Item(2).Table
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Luiz Antonio Pereira Silva, Joaquim José Cintra Maia Honório"
